# Applies sexting intensity escalation text updates to the "ZansiJourney" sheet
# (column B) per the PPV0 phase-detection bug fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ZansiJourney")

$ws.Range("B4").Value  = "cum with me daddy... I'm done holding back 🔥"
$ws.Range("B6").Value  = "I'm cumming for you right now daddy... FUCK watch this"
$ws.Range("B7").Value  = "I'm about to cum and I need you right here watching me when it happens 🔥"
$ws.Range("B8").Value  = "oh fuck 😏"
$ws.Range("B9").Value  = "you're not ready for this but I'm showing you anyway 🔥"
$ws.Range("B11").Value = "I can feel it coming and I'm not holding back... you need to see this"
$ws.Range("B12").Value = "I want you so deep inside me I can feel it in my chest... god I'm going crazy 🔥"
$ws.Range("B13").Value = "I'm fucking myself right now and all I can think about is you watching me do it daddy"
$ws.Range("B14").Value = "jesus fuck 😏"
$ws.Range("B15").Value = "this is what your words do to me daddy... watch 🔥"
$ws.Range("B17").Value = "tell me what you want me to do next... be specific, I want to hear every word"
$ws.Range("B18").Value = "I'm soaking wet and my fingers are going in and out and it's not enough daddy... I need you 🔥"
$ws.Range("B19").Value = "I literally can't stop now... I'm so turned on my whole body is aching for it"
$ws.Range("B20").Value = "mm okay wow... that hit different 😏"
$ws.Range("B21").Value = "you asked for more daddy... be careful what you wish for 🔥"
$ws.Range("B23").Value = "I'm sliding my hand between my legs right now and I'm already wet for you daddy"
$ws.Range("B24").Value = "the way you reacted... it's making me feel things all over my body right now 😏"
$ws.Range("B25").Value = "and? I can already tell you want more 🔥"
